$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns so numeric-looking
# strings (e.g. "6.00", "59.850.41") are not auto-converted to numbers.
$priceRange = $ws.Range('D2:D51')
$priceRange.NumberFormat = '@'

$ws.Range('D2').Value = '59.850.41'
$ws.Range('E2').Value = '  +1.62%  '
$ws.Range('D3').Value = '2.551.09'
$ws.Range('E3').Value = '  +3.38%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '502.19'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').Value = '151.44'
$ws.Range('E6').Value = '  -5.23%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('E8').Value = '  -6.57%  '
$ws.Range('D9').Value = '2.562.64'
$ws.Range('E9').Value = '  +2.81%  '
$ws.Range('D10').Value = '6.73'
$ws.Range('E10').Value = '  +6.91%  '
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('D12').Value = '0.341'
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '3.006.78'
$ws.Range('E14').Value = '  +4.19%  '
$ws.Range('D15').Value = '59.959.57'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('D16').Value = '21.38'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').Value = '2.565.55'
$ws.Range('E18').Value = '  +3.15%  '
$ws.Range('D19').Value = '4.77'
$ws.Range('E19').Value = '  +0.65%  '
$ws.Range('D20').Value = '344.54'
$ws.Range('E20').Value = '  +4.18%  '
$ws.Range('D21').Value = '10.18'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '6.00'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').Value = '59.77'
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('E25').Value = '  +0.85%  '
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('D27').Value = '2.683.85'
$ws.Range('E27').Value = '  +4.48%  '
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').Value = '0.0₃0838'
$ws.Range('E29').Value = '  +3.42%  '
$ws.Range('D30').Value = '7.39'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('D32').Value = '154.79'
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('D33').Value = '19.06'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('D34').Value = '1.54'
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D35').Value = '5.67'
$ws.Range('E35').Value = '  +3.54%  '
$ws.Range('D36').Value = '3.95'
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = '0.847'
$ws.Range('E38').Value = '  +18.70%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').Value = '3.73'
$ws.Range('E39').Value = '  +1.78%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '0.835'
$ws.Range('E40').Value = '  -1.91%  '
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('E42').Value = '  +2.42%  '
$ws.Range('D43').Value = '295.12'
$ws.Range('E43').Value = '  +3.19%  '
$ws.Range('D44').Value = '0.0558'
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '0.0989'
$ws.Range('E46').Value = '  -2.29%  '
$ws.Range('D47').Value = '0.611'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('D48').Value = '19.48'
$ws.Range('E48').Value = '  +6.99%  '
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('D50').Value = '0.0232'
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('D51').Value = '10.26'
$ws.Range('E51').Value = '  -0.12%  '

# Restore default (unstyled) cell formatting to match original style indices
$priceRange.Style = 'Normal'
